# Auto-generated edit script: updates cryptos list (price & 1h volume change)
# columns per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.247.16"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "1.574.84"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'207.41"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'22.36"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'0.0864"
$ws.Range("D12").Value = "1.799.39"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "1.567.93"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "'62.63"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "27.260.18"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").Value = "'216.04"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'9.42"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "'151.86"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "'6.68"
$ws.Range("E26").Value = "  -6.65%  "
$ws.Range("D27").Value = "'14.96"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").Value = "1.406.10"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "'0.943"
$ws.Range("E36").Value = "  -2.92%  "
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'0.519"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").Value = "'5.35"
$ws.Range("E44").Value = "  +1.70%  "
$ws.Range("D45").Value = "'63.94"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "1.711.27"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "'86.25"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").Value = "0.0₇0975"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").Value = "'0.0954"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -0.44%  "
